$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: add new hyperlink cell C15 -> matrix-diagonal-sum
$ws.Range("C15").Value = "https://leetcode.com/problems/matrix-diagonal-sum/"
$ws.Hyperlinks.Add($ws.Range("C15"), "https://leetcode.com/problems/matrix-diagonal-sum/")
$ws.Range("C15").Style = "Hyperlink"
$ws.Rows.Item(15).RowHeight = 15.75

# Row 16: add new hyperlink cell C16 -> transpose-matrix (A16 already has content)
$ws.Range("C16").Value = "https://leetcode.com/problems/transpose-matrix/"
$ws.Hyperlinks.Add($ws.Range("C16"), "https://leetcode.com/problems/transpose-matrix/")
$ws.Range("C16").Style = "Hyperlink"

# Update selection to D14
$ws.Range("D14").Select()
